# Fix Training Data Issue (#48)
# The BF column ("Date") stored the source filename-derived label
# "4-24-2007-08" for every data row. NBA.com's stats page showed the
# games under the day *before* they actually happened (because of how
# the box scores were posted), so the correct calendar date for this
# sheet is 2008-04-24. Re-stamp every data row (BF2:BF31) with the
# corrected date, written as literal text (not an Excel date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctedDate = "2008-04-24"
$dataRange = $ws.Range("BF2:BF31")

# Force text storage first -- otherwise Excel's normal "typed value"
# auto-detection would silently reinterpret an unformatted
# "2008-04-24" entry as a date serial number instead of keeping it as
# the plain string the source data uses. Once the literal text is in
# place, drop the temporary number format again so the cells keep
# their original (default) style, exactly like every other cell in
# this column.
$dataRange.NumberFormat = "@"
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = $correctedDate
}
$dataRange.ClearFormats()
